$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row values (B1:E1)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Update row 2 values (B2:E2)
$ws.Range("B2").Value = 7.9019510033225515
$ws.Range("C2").Value = 2.0966754236930001
$ws.Range("D2").Value = 1.983197291008338
$ws.Range("E2").Value = 1.933826288135009

# Update row 3 values (B3:E3)
$ws.Range("B3").Value = 6.8979008167237339
$ws.Range("C3").Value = 13.079370558231119
$ws.Range("D3").Value = 12.786249941062096
$ws.Range("E3").Value = -3.4107221288695655

# Update the selection to match the new range
$ws.Range("B1:E3").Select()
